# "Remove Data Row Activity"
# The original "Sheet1" master table (Roll No./Name/Dept/Sem1/Sem2/Sem3/Avg)
# had its row for Roll No. 1021 (Kumatsu / Current Affairs) removed, and the
# resulting data table was written out to a brand-new worksheet ("Sheet4"),
# inserted right before "Sheet1". This mirrors UiPath's "Remove Data Row" +
# "Write Range" activities: the new sheet gets plain literal values (no
# styles, no formulas, default column widths) and becomes the active sheet.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")

# Insert the new sheet right before "Sheet1" so the tab order becomes
# Sheet2, Sheet3, Sheet4, Sheet1.
$newSheet = $wb.Worksheets.Add($sheet1)
$newSheet.Name = "Sheet4"

$headers = @("Roll No.", "Name", "Dept", "Sem 1", "Sem 2", "Sem 3", "Avg")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Original 14-row student table minus Roll No. 1021 (Kumatsu).
$data = @(
    @(1011, "Ali",      "Micro Biology",   18, 18, 18, 18),
    @(1018, "Boris",    "Politics",        10, 15, 19, 14.666666666666666),
    @(1015, "Dialo",    "Mechanical",      18, 15, 18, 17),
    @(1020, "Emanual",  "Politics",         2,  5,  4, 3.6666666666666665),
    @(1022, "Fatumata", "Socialogy",       12, 18, 20, 16.666666666666668),
    @(1010, "Premji",   "Business",        19, 19, 20, 19.333333333333332),
    @(1019, "Salman",   "Politics",        12, 13, 15, 13.333333333333334),
    @(1016, "Sameera",  "Micro Biology",   17, 18, 19, 18),
    @(1012, "Selvi",    "Home Science",    15, 17, 18, 16.666666666666668),
    @(1013, "Vajpay",   "Language",        14, 13, 19, 15.333333333333334),
    @(1014, "Yi",       "IT",              17, 15, 18, 16.666666666666668),
    @(1023, "Yi Wong",  "Current Affairs", 17, 15, 18, 16.666666666666668),
    @(1017, "Ying",     "Language",        15, 10, 18, 14.333333333333334)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newSheet.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# Matches the saved selection on the new active sheet.
$newSheet.Range("N20").Select()
